# Collapse the word-by-word runs in the Title, Author and Abstract
# paragraphs into a single run each, matching the target rendering.
#
# Each of these paragraphs currently stores its text as a sequence of
# many <w:r> runs (one per word/space) that all share identical - i.e.
# empty/default - run formatting. We merge them into one run while
# keeping the *first* run element untouched (so its existing
# xml:space="preserve" attribute on <w:t> survives), deleting the
# remaining runs, and appending the rest of the text onto the end of
# that first run with InsertAfter.

function Merge-ParagraphRuns($para, $firstRunLength, $fullText) {
    $paraRange = $para.Range
    $paraRange.MoveEnd(1, -1)  # wdCharacter: drop the trailing paragraph mark

    # Range covering everything after the first run - delete it, leaving
    # just the original first run (formatting + its xml:space) intact.
    $tail = $paraRange.Duplicate
    $tail.Start = $paraRange.Start + $firstRunLength
    if ($tail.Start -lt $tail.End) {
        $tail.Delete()
    }

    # Re-grab the (now shorter) paragraph range and tack the remainder
    # of the target text onto the surviving first run.
    $remaining = $fullText.Substring($firstRunLength)
    if ($remaining.Length -gt 0) {
        $headRange = $para.Range
        $headRange.MoveEnd(1, -1)
        $headRange.InsertAfter($remaining)
    }
}

$d = $word.ActiveDocument

$title = $d.Paragraphs(1)
Merge-ParagraphRuns $title 8 "Answers: Laws of indices"

$author = $d.Paragraphs(2)
Merge-ParagraphRuns $author 8 "Isabella Lewis, Akshat Srivastava"

$abstract = $d.Paragraphs(4)
Merge-ParagraphRuns $abstract 7 "Answers to questions relating to using laws of indices."
